# Apply the ECOG performance status ValueSet edit to the workbook.
# (Mirrors an automated FHIR Implementation Guide re-publish of the
#  "Metadata" / "Include from LOINC" value-set spreadsheet.)

$wb = $excel.ActiveWorkbook

$wsMeta    = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from LOINC")

# --- Metadata sheet updates --------------------------------------------------
# URL (row 2)
$wsMeta.Range("B2").Value = "http://idg-rlp.de/fhir/tumorkonferenzen/ValueSet/ecog-performance-status"

# Date (row 8)
$wsMeta.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new "Jurisdiction" row just before "Description" (currently row 11),
# copying the formatting of the row above so the new row matches the sheet's
# existing body-row style.
$wsMeta.Rows("11").Insert()
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$wsMeta.Range("A11").Value = "Jurisdiction"
# The Jurisdiction value is an (explicit) empty string, not a blank cell -
# grab one from the sheet (an existing empty-valued cell, row 8 of the
# Include sheet, before it is touched below) so it round-trips as a real
# empty string rather than a cleared cell.
$wsInclude.Range("A8").Copy()
$wsMeta.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Include from LOINC sheet updates ---------------------------------------
# Insert a new concept row ("LA4489-6" / "ECOG Status unbekannt") before the
# existing blank row (currently row 8).
$wsInclude.Rows("8").Insert()
$wsInclude.Range("A7:B7").Copy()
$wsInclude.Range("A8:B8").PasteSpecial(-4122)  # xlPasteFormats
$wsInclude.Range("A8").Value = "LA4489-6"
$wsInclude.Range("B8").Value = "ECOG Status unbekannt"
$excel.CutCopyMode = 0

# --- Rename the "Include from LOINC" sheet to "Include #0" -----------------
$wsInclude.Name = "Include #0"
